$wb = $excel.ActiveWorkbook

# --- Update effort values (hours) on "Casos de Uso" sheet ---
$ws1 = $wb.Worksheets.Item("Casos de Uso")

$ws1.Range("F14").Value = 1.1000000000000001
$ws1.Range("F15").Value = 0.57999999999999996
$ws1.Range("F16").Value = 0.58599999999999997
$ws1.Range("F17").Value = 1.25

# --- Update the view / selection state to match where the author left off ---
$ws1.Activate()
$ws1.Range("C18").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
